$d = $word.ActiveDocument

# Locate the "Level 6" paragraph (unique in the document) - the three
# paragraphs we need to replace immediately follow it:
#   1) #{docu[bookmark]ment_a:document}
#   2) #{document_b:document}
#   3) an empty (it-IT) paragraph
$anchor = $d.Content
$anchor.Find.Execute("Level 6") | Out-Null
$anchor.Expand(4) | Out-Null  # wdParagraph

$p1 = $d.Range($anchor.End, $anchor.End)
$p1.Expand(4) | Out-Null
$p2 = $d.Range($p1.End, $p1.End)
$p2.Expand(4) | Out-Null
$p3 = $d.Range($p2.End, $p2.End)
$p3.Expand(4) | Out-Null

$target = $d.Range($p1.Start, $p3.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:ind w:left="1440"/>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>#{</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>document_</w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>a:document</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:ind w:left="1440"/>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>#{</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>document_</w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>b:document</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>#{</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>list:asunique</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Paragrafoelenco"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>#{</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>list:asnotunique</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr><w:lang w:val="en-US"/></w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)
